$d = $word.ActiveDocument

# -----------------------------------------------------------------------
# Update resume location: "Manchester, CT" -> "Alliance, OH" and relocate
# the stray "_GoBack" bookmark from the Skills line to the address line
# (matching the upstream commit "update to resume location").
# -----------------------------------------------------------------------

# Step 1: type "Alliance, " right after the existing "Manchester, CT" text.
# Inserting here (inside the run, not at its left edge) lets the new text
# inherit the correct Roboto/20 run formatting instead of bleeding into
# the neighboring <w:br/> run.
$r = $d.Content
$r.Find.Execute("Manchester, CT", $true, $false, $false, $false, $false, `
                 $true, 1, $false, "", 0)
$r.Collapse(0)
$r.InsertBefore("Alliance, ")

# Step 2: remove the old "Manchester, CT" text (plain deletion, no new
# text typed at the run boundary, so formatting of neighbors is untouched).
$r2 = $d.Content
$r2.Find.Execute("Manchester, CT", $true, $false, $false, $false, $false, `
                  $true, 1, $false, "", 1)

# Step 3: relocate the "_GoBack" bookmark so it sits right after
# "Alliance, " (it currently wraps the end of "Unity 3D, Eclipse").
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

$r3 = $d.Content
$r3.Find.Execute("Alliance, ", $true, $false, $false, $false, $false, `
                  $true, 1, $false, "", 0)
$r3.Collapse(0)
$bmStart = $r3.Start
$d.Bookmarks.Add("_GoBack", $r3)

# Step 4: type "OH | " right before "518-334-8045" (again a safe, interior
# insertion point so it inherits the surrounding Roboto/20 formatting).
$r4 = $d.Content
$r4.Find.Execute("518-334-8045", $true, $false, $false, $false, $false, `
                  $true, 1, $false, "", 0)
$numStart = $r4.Start
$r4.Collapse(1)
$r4.InsertBefore("OH | ")

# Step 5: delete the old " | " separator that used to follow
# "Manchester, CT" (now sitting between the bookmark and the freshly
# typed "OH | ").
$oldSep = $d.Range($bmStart, $numStart)
$oldSep.Delete()

# Step 6: split "OH | 518-334-8045 | " into discrete runs "OH ", "|",
# " " (and leave "518-334-8045 | " as the trailing run) by toggling a
# character property on/off, which forces a run break without altering
# the visible formatting.
$r6 = $d.Content
$r6.Find.Execute("OH | 518", $true, $false, $false, $false, $false, `
                  $true, 1, $false, "", 0)
$ohStart = $r6.Start

$p1 = $d.Range($ohStart, $ohStart + 3)
$p1.Bold = 1
$p1.Bold = 0

$p2 = $d.Range($ohStart + 3, $ohStart + 4)
$p2.Bold = 1
$p2.Bold = 0

$p3 = $d.Range($ohStart + 4, $ohStart + 5)
$p3.Bold = 1
$p3.Bold = 0
